$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D7").Value = -7.731600000000001
$ws.Range("A9").Value = -21.9456
$ws.Range("D12").Value = -7.180200000000002
$ws.Range("E15").Value = 16.26540000000001
$ws.Range("A18").Value = -22.35770000000001
$ws.Range("A20").Value = -19.6457
$ws.Range("D26").Value = -8.667299999999999
$ws.Range("A27").Value = -22.0723
$ws.Range("D27").Value = -8.853199999999999
$ws.Range("D29").Value = -7.241400000000001
$ws.Range("D37").Value = -7.6882
$ws.Range("D38").Value = -7.365600000000001
$ws.Range("E38").Value = 16.72639999999999
$ws.Range("E44").Value = 16.3672
$ws.Range("D51").Value = -8.695499999999999
$ws.Range("E51").Value = 16.28939999999999
$ws.Range("D55").Value = -8.8581
$ws.Range("E57").Value = 16.7064
$ws.Range("E63").Value = 18.41470000000001
$ws.Range("A69").Value = -21.5872
$ws.Range("D69").Value = -7.155399999999994
$ws.Range("D70").Value = -7.565700000000003
$ws.Range("E70").Value = 17.07540000000001
$ws.Range("A76").Value = -19.88249999999999
$ws.Range("A82").Value = -21.9426
$ws.Range("D83").Value = -9.0062
$ws.Range("E99").Value = 16.5182
$ws.Range("D102").Value = -7.7389
